# Fix the typo'd category label in A22 ("Residensialtial" -> "residential"),
# which is the substantive content change behind this commit; the shared
# string table reindexing that shows up across the rest of the sheet is a
# mechanical side effect of that single edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A22").Value = "residential"

# Widen column A to fit the edited text and record the cell the user left
# selected, same as Excel would persist after an interactive edit.
$ws.Columns("A").ColumnWidth = 17.333333333333332

$ws.Range("A22").Select()

# Reflect the scrolled viewport (topLeftCell = A10) on the active window.
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
